# Apply "Correction type pour génération à partir fsh" changes:
#  - Metadata!B4 (Name row) gets the value "CiviliteVs"
#  - Metadata!B8 (Date row) is updated to the new generation timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B4").Value = "CiviliteVs"
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
